$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing data to make room for headers
$ws.Rows.Item(1).Insert()

# Set the new header row values
$ws.Range("A1").Value = "Folder Name"
$ws.Range("B1").Value = "Keyword"

# Update the selection to match the target state
$ws.Range("F2").Select()
